# Generate Report for Handoff
# Adds two new handed-off files (5196cec1-... and c700ca98-...) to the
# Overview / zh-cn / de-de report sheets, growing each table from 4 data
# rows (A1:x5) to 6 data rows (A1:x7).

$wb = $excel.ActiveWorkbook
$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview"  (table3 -> rows 6 & 7, columns A:G)
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("G6").NumberFormat = $dateFmt
$ov.Range("G7").NumberFormat = $dateFmt

$ov.Range("A6").Value = "5196cec1-bec1-46cc-a8dd-122c0812e5f1.md"
$ov.Range("C6").Value = ".md"
$ov.Range("E6").Value = "Ready for handoff"
$ov.Range("F6").Value = "Ready for handoff"
$ov.Range("G6").Value = "2016-08-26 14:52:27"
$ov.Hyperlinks.Add($ov.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5196cec1bec146cca8dd122c0812e5f1000000/e2e/5196cec1-bec1-46cc-a8dd-122c0812e5f1.md", "", "", "e2e\5196cec1-bec1-46cc-a8dd-122c0812e5f1.md")

$ov.Range("A7").Value = "c700ca98-f34d-40d1-8ec8-05c12b21372e.md"
$ov.Range("C7").Value = ".md"
$ov.Range("E7").Value = "Ready for handoff"
$ov.Range("F7").Value = "Ready for handoff"
$ov.Range("G7").Value = "2016-08-26 14:52:27"
$ov.Hyperlinks.Add($ov.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c700ca98f34d40d18ec805c12b21372e0000000/e2e/c700ca98-f34d-40d1-8ec8-05c12b21372e.md", "", "", "e2e\c700ca98-f34d-40d1-8ec8-05c12b21372e.md")

$ovTable = $ov.ListObjects.Item(1)
$ovTable.Resize($ov.Range("A1:G7"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"  (table1 -> rows 6 & 7, columns A:P)
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("H6").NumberFormat = $dateFmt
$zh.Range("H7").NumberFormat = $dateFmt
$zh.Range("K6").NumberFormat = $dateFmt
$zh.Range("K7").NumberFormat = $dateFmt

$zh.Range("B6").Value = ".md"
$zh.Range("C6").Value = "Ready for handoff"
$zh.Range("D6").Value = "e2e"
$zh.Range("E6").Value = "ht"
$zh.Range("F6").Value = "False"
$zh.Range("G6").Value = "5196cec1-bec1-46cc-a8dd-122c0812e5f1.7caec08c899470a7d086f11d61a7efb2af0da9cc.zh-cn.xlf"
$zh.Range("H6").Value = "2016-08-26 14:52:23"
$zh.Range("K6").Value = "0001-01-01 00:00:00"
$zh.Range("M6").Value = "True"
$zh.Range("O6").Value = "False"
$zh.Hyperlinks.Add($zh.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5196cec1bec146cca8dd122c0812e5f1000000/e2e/5196cec1-bec1-46cc-a8dd-122c0812e5f1.md", "", "", "5196cec1-bec1-46cc-a8dd-122c0812e5f1.md")

$zh.Range("B7").Value = ".md"
$zh.Range("C7").Value = "Ready for handoff"
$zh.Range("D7").Value = "e2e"
$zh.Range("E7").Value = "ht"
$zh.Range("F7").Value = "False"
$zh.Range("G7").Value = "c700ca98-f34d-40d1-8ec8-05c12b21372e.6af597acc4103b35879ae6172741ea614c06d3be.zh-cn.xlf"
$zh.Range("H7").Value = "2016-08-26 14:52:23"
$zh.Range("K7").Value = "0001-01-01 00:00:00"
$zh.Range("M7").Value = "True"
$zh.Range("O7").Value = "False"
$zh.Hyperlinks.Add($zh.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c700ca98f34d40d18ec805c12b21372e0000000/e2e/c700ca98-f34d-40d1-8ec8-05c12b21372e.md", "", "", "c700ca98-f34d-40d1-8ec8-05c12b21372e.md")

$zhTable = $zh.ListObjects.Item(1)
$zhTable.Resize($zh.Range("A1:P7"))

# ---------------------------------------------------------------------
# Sheet "de-de"  (table2 -> rows 6 & 7, columns A:P)
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("H6").NumberFormat = $dateFmt
$de.Range("H7").NumberFormat = $dateFmt
$de.Range("K6").NumberFormat = $dateFmt
$de.Range("K7").NumberFormat = $dateFmt

$de.Range("B6").Value = ".md"
$de.Range("C6").Value = "Ready for handoff"
$de.Range("D6").Value = "e2e"
$de.Range("E6").Value = "ht"
$de.Range("F6").Value = "False"
$de.Range("G6").Value = "5196cec1-bec1-46cc-a8dd-122c0812e5f1.7caec08c899470a7d086f11d61a7efb2af0da9cc.de-de.xlf"
$de.Range("H6").Value = "2016-08-26 14:52:27"
$de.Range("K6").Value = "0001-01-01 00:00:00"
$de.Range("M6").Value = "True"
$de.Range("O6").Value = "False"
$de.Hyperlinks.Add($de.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5196cec1bec146cca8dd122c0812e5f1000000/e2e/5196cec1-bec1-46cc-a8dd-122c0812e5f1.md", "", "", "5196cec1-bec1-46cc-a8dd-122c0812e5f1.md")

$de.Range("B7").Value = ".md"
$de.Range("C7").Value = "Ready for handoff"
$de.Range("D7").Value = "e2e"
$de.Range("E7").Value = "ht"
$de.Range("F7").Value = "False"
$de.Range("G7").Value = "c700ca98-f34d-40d1-8ec8-05c12b21372e.6af597acc4103b35879ae6172741ea614c06d3be.de-de.xlf"
$de.Range("H7").Value = "2016-08-26 14:52:27"
$de.Range("K7").Value = "0001-01-01 00:00:00"
$de.Range("M7").Value = "True"
$de.Range("O7").Value = "False"
$de.Hyperlinks.Add($de.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c700ca98f34d40d18ec805c12b21372e0000000/e2e/c700ca98-f34d-40d1-8ec8-05c12b21372e.md", "", "", "c700ca98-f34d-40d1-8ec8-05c12b21372e.md")

$deTable = $de.ListObjects.Item(1)
$deTable.Resize($de.Range("A1:P7"))
